# Staff list: add header row, restructure columns, make the header bold,
# turn the range into an AutoFilter, and sort the data (A-Z) on the
# "LAST NAME" column - see commit message: "added sorting to excel file
# for when creating the staff list will automatically be in alphabetical
# order".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------
$ws.Range("A1").Value = "FIRST NAME"
$ws.Range("B1").Value = "LAST NAME"
$ws.Range("C1").Value = "CODE"
$ws.Range("D1").Value = "ROLE"

# --- Data rows --------------------------------------------------------
# Row 2 becomes the ADMIN / 1111 account (replacing Nathan Danskin's row).
$ws.Range("A2").Value = "ADMIN"
$ws.Range("B2").Value = "ADMIN"
$ws.Range("C2").Value = "1111"
$ws.Range("D2").Value = "ADMIN"

# New rows with just a single value in column B, as per the source data.
$ws.Range("B3").Value = "dead"
$ws.Range("B4").Value = "roll"

# --- Header formatting (bold) -----------------------------------------
$ws.Range("A1:D1").Font.Bold = $true

# --- Column widths (approximate autofit-style resize) ------------------
$ws.Columns.Item(1).ColumnWidth = 10.88
$ws.Columns.Item(2).ColumnWidth = 10.16
$ws.Columns.Item(3).ColumnWidth = 5.02
$ws.Columns.Item(4).ColumnWidth = 6.59

# --- AutoFilter + sort --------------------------------------------------
$rng = $ws.Range("A1:D4")
$rng.AutoFilter()

# Sort ascending by column B (LAST NAME), header row excluded from the sort.
$rng.Sort($ws.Range("B1:B4"), 1, $null, $null, 1, 0, 0, 1)

# Record the (hidden) _FilterDatabase defined name the way Excel does
# whenever AutoFilter is applied, scoped to this sheet.
$fd = $ws.Names.Add("_xlnm._FilterDatabase", "=Staff!`$A`$1:`$D`$4")
$fd.Visible = $false

Write-Host "done"
